# Oxo_Dimer_Energies.xlsx edit:
#  - Change energy unit label from kcal/mol to kJ/mol
#  - Update the conversion-factor formulas in column G accordingly
#  - Refresh the sheet view (zoom level + selected cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "Relative Energy" header (column G, row 1) from kcal/mol to kJ/mol
$ws.Range("G1").Value = "Relative Energy / kJ/mol"

# Update the relative-energy formulas to use the kJ/mol conversion factor
# (1 Hartree = 2625.5 kJ/mol, previously 627.5095 kcal/mol)
$ws.Range("G5").Formula = "=((D5+D4)-2*D2)*2625.5"
$ws.Range("G6").Formula = "=((D5+D4)-D3)*2625.5"

# Update the sheet view: zoom in to 125% and move the active selection
$excel.ActiveWindow.Zoom = 125
$ws.Range("G7").Select() | Out-Null
